$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 147, pushing existing rows 147:207 down to 149:209.
$ws.Rows("147:148").Insert()

# New row 147: Primera, week of 2021-09-27 (serial 44466)
$ws.Cells.Item(147, 1).Value = 1
$ws.Cells.Item(147, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(147, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(147, 4).Value = 44466
$ws.Cells.Item(147, 5).Value = 15
$ws.Cells.Item(147, 6).Value = 100112032
$ws.Cells.Item(147, 7).Value = "Zapallo italiano"
$ws.Cells.Item(147, 8).Value = "Huracán"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 130
$ws.Cells.Item(147, 11).Value = 9000
$ws.Cells.Item(147, 12).Value = 10000
$ws.Cells.Item(147, 13).Value = 9500
$ws.Cells.Item(147, 14).Value = '$/caja 70 unidades'
$ws.Cells.Item(147, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(147, 16).Value = 136
$ws.Cells.Item(147, 17).Value = 70
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# New row 148: Segunda, week of 2021-09-27 (serial 44466)
$ws.Cells.Item(148, 1).Value = 1
$ws.Cells.Item(148, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(148, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(148, 4).Value = 44466
$ws.Cells.Item(148, 5).Value = 15
$ws.Cells.Item(148, 6).Value = 100112032
$ws.Cells.Item(148, 7).Value = "Zapallo italiano"
$ws.Cells.Item(148, 8).Value = "Huracán"
$ws.Cells.Item(148, 9).Value = "Segunda"
$ws.Cells.Item(148, 10).Value = 140
$ws.Cells.Item(148, 11).Value = 7000
$ws.Cells.Item(148, 12).Value = 8000
$ws.Cells.Item(148, 13).Value = 7500
$ws.Cells.Item(148, 14).Value = '$/caja 100 unidades'
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 75
$ws.Cells.Item(148, 17).Value = 100
$ws.Cells.Item(148, 18).Value = "Hortaliza"
